# repull data, push all data, mean calculation
# Updates the "dSF" column (F) with recalculated values for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dsfValues = @{
    2  = -4
    3  = 4
    4  = 3
    5  = 1
    6  = -1
    7  = 3
    8  = 2
    9  = -1
    10 = -1
    11 = -1
    12 = -2
    13 = 1
    14 = -2
    15 = 2
    16 = 7
    17 = -1
    18 = -4
    19 = -1
    20 = 6
    21 = -3
    22 = 9
    23 = 2
    24 = 4
    25 = -1
    26 = 2
    27 = -3
    28 = 0
    30 = 1
    31 = 0
    32 = 2
    33 = 2
    35 = 1
    36 = 3
    37 = -1
}

foreach ($row in $dsfValues.Keys) {
    $ws.Range("F$row").Value = $dsfValues[$row]
}
